$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Fail Log message for row 13 (selector #lblServiceID, chrome-headless-shell 122.0.6261.129)
$newFailLog13 = "no such element: Unable to locate element: {`"method`":`"css selector`",`"selector`":`"#lblServiceID`"}`n  (Session info: chrome-headless-shell=122.0.6261.129)`nFor documentation on this error, please visit: https://www.seleniumhq.org/exceptions/no_such_element.html`nBuild info: version: '3.141.59', revision: 'e82be7d358', time: '2018-11-14T08:17:03'`nSystem info: host: 'SIPL92', ip: '10.212.130.91', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '20'`nDriver info: org.openqa.selenium.chrome.ChromeDriver`nCapabilities {acceptInsecureCerts: false, browserName: chrome-headless-shell, browserVersion: 122.0.6261.129, chrome: {chromedriverVersion: 122.0.6261.128 (f18a44fedeb..., userDataDir: C:\Users\RPRAJA~1\AppData\L...}, fedcm:accounts: true, goog:chromeOptions: {debuggerAddress: localhost:59519}, javascriptEnabled: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: WINDOWS, platformName: WINDOWS, proxy: Proxy(), setWindowRect: true, strictFileInteractability: false, timeouts: {implicit: 0, pageLoad: 300000, script: 30000}, unhandledPromptBehavior: dismiss and notify, webauthn:extension:credBlob: true, webauthn:extension:largeBlob: true, webauthn:extension:minPinLength: true, webauthn:extension:prf: true, webauthn:virtualAuthenticators: true}`nSession ID: 20ab763d4621b4c19b079e63bf171860`n*** Element info: {Using=id, value=lblServiceID}"

# Row 13: PickupID/POD No changes (stored as text, like the rest of the column), Fail Log changes to the new selenium error
# Use a quoted-text formula + paste-as-values so the result lands as a genuine
# text cell (matching the rest of column C) without leaving the cell's style
# touched the way a NumberFormat round-trip would.
$ws.Cells.Item(13, 3).Formula = '="15287848"'
$ws.Cells.Item(13, 3).Copy() | Out-Null
$ws.Cells.Item(13, 3).PasteSpecial(-4163) | Out-Null
$ws.Range("F13").Value = $newFailLog13
$ws.Rows(13).AutoFit()

# Row 14: PickupID/POD No changes, Connect OP Result flips to FAIL, Fail Log becomes a NullPointer message
$ws.Cells.Item(14, 3).Formula = '="15287826"'
$ws.Cells.Item(14, 3).Copy() | Out-Null
$ws.Cells.Item(14, 3).PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "FAIL"
$ws.Range("F14").Value = 'Cannot invoke "org.openqa.selenium.WebElement.isDisplayed()" because "element" is null'
$ws.Rows(14).AutoFit()

$excel.CutCopyMode = $false
